$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New revision history rows (B column values are serials for
# 2012-11-22 / 2012-11-26 / 2012-11-27, formatted as dates below)
$ws.Range("B3").Value = 41235
$ws.Range("C3").Value = "Added use case diagrams, definition of done, vision and glossery"

$ws.Range("B4").Value = 41239
$ws.Range("C4").Value = "Added versioning conversions and use cases"

$ws.Range("B5").Value = 41240
$ws.Range("C5").Value = "Added product backlog, changed use cases, added estimates for tasks and sprint backlog"

# Apply date number format to the first new cell, then propagate it (format only)
# to the rest so all three share a single style entry (numFmtId 14 - built in
# short date format), matching how Excel dedupes identical cell styles.
$ws.Range("B3").NumberFormat = "mm-dd-yy"
$ws.Range("B3").Copy()
$ws.Range("B4:B5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Widen the description column to fit the new text (closest value the
# pixel-grid column-width conversion can reach to the authored 79.140625)
$ws.Columns.Item(3).ColumnWidth = 78.33

# Move the selection like the author's last interaction
$null = $ws.Range("F6").Select()
